$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: seed the two new columns (L, M) with the formatting of column K
# (column K is the last pre-existing data column) across the three
# financial-statement blocks, so the new cells inherit the correct style.
$ws.Range("K7:K35").Copy($ws.Range("L7:L35"))
$ws.Range("K38:K77").Copy($ws.Range("L38:L77"))
$ws.Range("K80:K102").Copy($ws.Range("L80:L102"))
$ws.Range("K7:K35").Copy($ws.Range("M7:M35"))
$ws.Range("K38:K77").Copy($ws.Range("M38:M77"))
$ws.Range("K80:K102").Copy($ws.Range("M80:M102"))

# Step 2: write the refreshed quarterly figures (2 new quarters added,
# historical quarters shifted right, some restated) for D:M per row.
$row7 = New-Object 'object[,]' 1,10
$row7[0,0] = 43465
$row7[0,1] = 43373
$row7[0,2] = 43281
$row7[0,3] = 43190
$row7[0,4] = 43100
$row7[0,5] = 43008
$row7[0,6] = 42916
$row7[0,7] = 42825
$row7[0,8] = 42735
$row7[0,9] = 42643
$ws.Range("D7:M7").Value2 = $row7

$row8 = New-Object 'object[,]' 1,10
$row8[0,0] = 2821200
$row8[0,1] = 2495200
$row8[0,2] = 2563500
$row8[0,3] = 2353200
$row8[0,4] = 2483700
$row8[0,5] = 2370600
$row8[0,6] = 1230700
$row8[0,7] = 1182200
$row8[0,8] = 1084600
$row8[0,9] = 1079900
$ws.Range("D8:M8").Value2 = $row8

$row9 = New-Object 'object[,]' 1,10
$row9[0,0] = 2345800
$row9[0,1] = 1979800
$row9[0,2] = 2181900
$row9[0,3] = 2052100
$row9[0,4] = 2039600
$row9[0,5] = 2032700
$row9[0,6] = 1185500
$row9[0,7] = 1038600
$row9[0,8] = 1080700
$row9[0,9] = 964000
$ws.Range("D9:M9").Value2 = $row9

$row10 = New-Object 'object[,]' 1,10
$row10[0,0] = 475400
$row10[0,1] = 515400
$row10[0,2] = 381600
$row10[0,3] = 301100
$row10[0,4] = 444100
$row10[0,5] = 337900
$row10[0,6] = 45200
$row10[0,7] = 143600
$row10[0,8] = 3900
$row10[0,9] = 115900
$ws.Range("D10:M10").Value2 = $row10

$row12 = New-Object 'object[,]' 1,10
$row12[0,0] = "NA"
$row12[0,1] = "NA"
$row12[0,2] = "NA"
$row12[0,3] = "NA"
$row12[0,4] = "NA"
$row12[0,5] = "NA"
$row12[0,6] = "NA"
$row12[0,7] = "NA"
$row12[0,8] = "NA"
$row12[0,9] = "NA"
$ws.Range("D12:M12").Value2 = $row12

$row13 = New-Object 'object[,]' 1,10
$row13[0,0] = 0
$row13[0,1] = 0
$row13[0,2] = 0
$row13[0,3] = 0
$row13[0,4] = 0
$row13[0,5] = 0
$row13[0,6] = 0
$row13[0,7] = 0
$row13[0,8] = 0
$row13[0,9] = 0
$ws.Range("D13:M13").Value2 = $row13

$row14 = New-Object 'object[,]' 1,10
$row14[0,0] = -100
$row14[0,1] = 100
$row14[0,2] = -13200
$row14[0,3] = 36500
$row14[0,4] = 0
$row14[0,5] = 0
$row14[0,6] = 0
$row14[0,7] = 0
$row14[0,8] = 0
$row14[0,9] = 245300
$ws.Range("D14:M14").Value2 = $row14

$row15 = New-Object 'object[,]' 1,10
$row15[0,0] = 53000
$row15[0,1] = 49200
$row15[0,2] = 49200
$row15[0,3] = 48000
$row15[0,4] = 47900
$row15[0,5] = 46900
$row15[0,6] = 29500
$row15[0,7] = 29000
$row15[0,8] = 29800
$row15[0,9] = 29000
$ws.Range("D15:M15").Value2 = $row15

$row17 = New-Object 'object[,]' 1,10
$row17[0,0] = 2638300
$row17[0,1] = 2240100
$row17[0,2] = 2415200
$row17[0,3] = 2350900
$row17[0,4] = 2371200
$row17[0,5] = 2286000
$row17[0,6] = 1277200
$row17[0,7] = 1152400
$row17[0,8] = 1128700
$row17[0,9] = 1328000
$ws.Range("D17:M17").Value2 = $row17

$row18 = New-Object 'object[,]' 1,10
$row18[0,0] = 182900
$row18[0,1] = 255100
$row18[0,2] = 148300
$row18[0,3] = 2300
$row18[0,4] = 112500
$row18[0,5] = 84600
$row18[0,6] = -46500
$row18[0,7] = 29800
$row18[0,8] = -44100
$row18[0,9] = -248100
$ws.Range("D18:M18").Value2 = $row18

$row20 = New-Object 'object[,]' 1,10
$row20[0,0] = 5000
$row20[0,1] = 12900
$row20[0,2] = 3500
$row20[0,3] = 1400
$row20[0,4] = 5000
$row20[0,5] = 201500
$row20[0,6] = 2200
$row20[0,7] = 4100
$row20[0,8] = -8900
$row20[0,9] = -5000
$ws.Range("D20:M20").Value2 = $row20

$row21 = New-Object 'object[,]' 1,10
$row21[0,0] = 240900
$row21[0,1] = 317200
$row21[0,2] = 201000
$row21[0,3] = 51700
$row21[0,4] = 165400
$row21[0,5] = 333000
$row21[0,6] = -14800
$row21[0,7] = 62900
$row21[0,8] = -23200
$row21[0,9] = -224100
$ws.Range("D21:M21").Value2 = $row21

$row22 = New-Object 'object[,]' 1,10
$row22[0,0] = 30700
$row22[0,1] = 31200
$row22[0,2] = 31500
$row22[0,3] = 32500
$row22[0,4] = 31300
$row22[0,5] = 34100
$row22[0,6] = 14900
$row22[0,7] = 13600
$row22[0,8] = 13700
$row22[0,9] = 13900
$ws.Range("D22:M22").Value2 = $row22

$row23 = New-Object 'object[,]' 1,10
$row23[0,0] = 157200
$row23[0,1] = 236800
$row23[0,2] = 120300
$row23[0,3] = -28800
$row23[0,4] = 86200
$row23[0,5] = 252000
$row23[0,6] = -59200
$row23[0,7] = 20300
$row23[0,8] = -66700
$row23[0,9] = -267000
$ws.Range("D23:M23").Value2 = $row23

$row24 = New-Object 'object[,]' 1,10
$row24[0,0] = 37800
$row24[0,1] = 51500
$row24[0,2] = 22800
$row24[0,3] = -9600
$row24[0,4] = 26200
$row24[0,5] = 133500
$row24[0,6] = -27000
$row24[0,7] = 5000
$row24[0,8] = -34700
$row24[0,9] = -103300
$ws.Range("D24:M24").Value2 = $row24

$row25 = New-Object 'object[,]' 1,10
$row25[0,0] = 0
$row25[0,1] = 0
$row25[0,2] = 0
$row25[0,3] = 0
$row25[0,4] = 0
$row25[0,5] = 0
$row25[0,6] = 0
$row25[0,7] = 0
$row25[0,8] = 0
$row25[0,9] = 0
$ws.Range("D25:M25").Value2 = $row25

$row26 = New-Object 'object[,]' 1,10
$row26[0,0] = 119400
$row26[0,1] = 185300
$row26[0,2] = 97500
$row26[0,3] = -19200
$row26[0,4] = 60000
$row26[0,5] = 118500
$row26[0,6] = -32200
$row26[0,7] = 15300
$row26[0,8] = -32000
$row26[0,9] = -163700
$ws.Range("D26:M26").Value2 = $row26

$row27 = New-Object 'object[,]' 1,10
$row27[0,0] = 113600
$row27[0,1] = 178800
$row27[0,2] = 89900
$row27[0,3] = -34100
$row27[0,4] = 46000
$row27[0,5] = 108500
$row27[0,6] = -37900
$row27[0,7] = 11200
$row27[0,8] = -36600
$row27[0,9] = -167700
$ws.Range("D27:M27").Value2 = $row27

$row28 = New-Object 'object[,]' 1,10
$row28[0,0] = 0
$row28[0,1] = 0
$row28[0,2] = 0
$row28[0,3] = 0
$row28[0,4] = 0
$row28[0,5] = 0
$row28[0,6] = 0
$row28[0,7] = 0
$row28[0,8] = 0
$row28[0,9] = 0
$ws.Range("D28:M28").Value2 = $row28

$row29 = New-Object 'object[,]' 1,10
$row29[0,0] = 2500
$row29[0,1] = 1000
$row29[0,2] = -10800
$row29[0,3] = -800
$row29[0,4] = 165100
$row29[0,5] = -4100
$row29[0,6] = 0
$row29[0,7] = 0
$row29[0,8] = 80800
$row29[0,9] = 6000
$ws.Range("D29:M29").Value2 = $row29

$row30 = New-Object 'object[,]' 1,10
$row30[0,0] = 0
$row30[0,1] = 0
$row30[0,2] = 0
$row30[0,3] = 0
$row30[0,4] = 0
$row30[0,5] = 0
$row30[0,6] = 0
$row30[0,7] = 0
$row30[0,8] = 0
$row30[0,9] = 0
$ws.Range("D30:M30").Value2 = $row30

$row31 = New-Object 'object[,]' 1,10
$row31[0,0] = 0
$row31[0,1] = 0
$row31[0,2] = 0
$row31[0,3] = 0
$row31[0,4] = 0
$row31[0,5] = 0
$row31[0,6] = 0
$row31[0,7] = 0
$row31[0,8] = 0
$row31[0,9] = 0
$ws.Range("D31:M31").Value2 = $row31

$row32 = New-Object 'object[,]' 1,10
$row32[0,0] = -5000
$row32[0,1] = -12900
$row32[0,2] = -3500
$row32[0,3] = -1400
$row32[0,4] = -5000
$row32[0,5] = -201500
$row32[0,6] = -2200
$row32[0,7] = -4100
$row32[0,8] = 8900
$row32[0,9] = 5000
$ws.Range("D32:M32").Value2 = $row32

$row33 = New-Object 'object[,]' 1,10
$row33[0,0] = 116100
$row33[0,1] = 179800
$row33[0,2] = 79100
$row33[0,3] = -34900
$row33[0,4] = 211100
$row33[0,5] = 104400
$row33[0,6] = -37900
$row33[0,7] = 11200
$row33[0,8] = 44200
$row33[0,9] = -161700
$ws.Range("D33:M33").Value2 = $row33

$row34 = New-Object 'object[,]' 1,10
$row34[0,0] = 0
$row34[0,1] = 0
$row34[0,2] = 0
$row34[0,3] = 0
$row34[0,4] = 0
$row34[0,5] = 0
$row34[0,6] = 0
$row34[0,7] = 0
$row34[0,8] = 0
$row34[0,9] = 0
$ws.Range("D34:M34").Value2 = $row34

$row35 = New-Object 'object[,]' 1,10
$row35[0,0] = 116100
$row35[0,1] = 179800
$row35[0,2] = 79100
$row35[0,3] = -34900
$row35[0,4] = 211100
$row35[0,5] = 104400
$row35[0,6] = -37900
$row35[0,7] = 11200
$row35[0,8] = 44200
$row35[0,9] = -161700
$ws.Range("D35:M35").Value2 = $row35

$row38 = New-Object 'object[,]' 1,10
$row38[0,0] = 43465
$row38[0,1] = 43373
$row38[0,2] = 43281
$row38[0,3] = 43190
$row38[0,4] = 43100
$row38[0,5] = 43008
$row38[0,6] = 42916
$row38[0,7] = 42825
$row38[0,8] = 42735
$row38[0,9] = 42643
$ws.Range("D38:M38").Value2 = $row38

$row41 = New-Object 'object[,]' 1,10
$row41[0,0] = 1079300
$row41[0,1] = 1109100
$row41[0,2] = 1132800
$row41[0,3] = 1018000
$row41[0,4] = 931800
$row41[0,5] = 831700
$row41[0,6] = 572300
$row41[0,7] = 591400
$row41[0,8] = 689200
$row41[0,9] = 315300
$ws.Range("D41:M41").Value2 = $row41

$row42 = New-Object 'object[,]' 1,10
$row42[0,0] = 0
$row42[0,1] = 0
$row42[0,2] = 0
$row42[0,3] = 0
$row42[0,4] = 0
$row42[0,5] = 0
$row42[0,6] = 0
$row42[0,7] = 0
$row42[0,8] = 0
$row42[0,9] = 0
$ws.Range("D42:M42").Value2 = $row42

$row43 = New-Object 'object[,]' 1,10
$row43[0,0] = 514400
$row43[0,1] = 740000
$row43[0,2] = 802100
$row43[0,3] = 644100
$row43[0,4] = 581700
$row43[0,5] = 495800
$row43[0,6] = 288000
$row43[0,7] = 327300
$row43[0,8] = 266000
$row43[0,9] = 197200
$ws.Range("D43:M43").Value2 = $row43

$row44 = New-Object 'object[,]' 1,10
$row44[0,0] = 690900
$row44[0,1] = 923900
$row44[0,2] = 881500
$row44[0,3] = 946000
$row44[0,4] = 808400
$row44[0,5] = 693500
$row44[0,6] = 377400
$row44[0,7] = 397400
$row44[0,8] = 392400
$row44[0,9] = 375100
$ws.Range("D44:M44").Value2 = $row44

$row45 = New-Object 'object[,]' 1,10
$row45[0,0] = 135700
$row45[0,1] = 88400
$row45[0,2] = 132900
$row45[0,3] = 248400
$row45[0,4] = 289900
$row45[0,5] = 249600
$row45[0,6] = 74000
$row45[0,7] = 51200
$row45[0,8] = 49300
$row45[0,9] = 528900
$ws.Range("D45:M45").Value2 = $row45

$row46 = New-Object 'object[,]' 1,10
$row46[0,0] = 2420300
$row46[0,1] = 2861400
$row46[0,2] = 2949300
$row46[0,3] = 2856500
$row46[0,4] = 2611800
$row46[0,5] = 2270600
$row46[0,6] = 1311700
$row46[0,7] = 1367300
$row46[0,8] = 1396900
$row46[0,9] = 1416500
$ws.Range("D46:M46").Value2 = $row46

$row47 = New-Object 'object[,]' 1,10
$row47[0,0] = 163000
$row47[0,1] = 166100
$row47[0,2] = 177200
$row47[0,3] = 184800
$row47[0,4] = 186500
$row47[0,5] = 189800
$row47[0,6] = 379600
$row47[0,7] = 390400
$row47[0,8] = 367500
$row47[0,9] = 373700
$ws.Range("D47:M47").Value2 = $row47

$row48 = New-Object 'object[,]' 1,10
$row48[0,0] = 2194900
$row48[0,1] = 2142400
$row48[0,2] = 2103600
$row48[0,3] = 2093500
$row48[0,4] = 2140800
$row48[0,5] = 2147700
$row48[0,6] = 1078600
$row48[0,7] = 1089800
$row48[0,8] = 1103300
$row48[0,9] = 1113500
$ws.Range("D48:M48").Value2 = $row48

$row49 = New-Object 'object[,]' 1,10
$row49[0,0] = 962200
$row49[0,1] = 963600
$row49[0,2] = 963900
$row49[0,3] = 926800
$row49[0,4] = 917700
$row49[0,5] = 888600
$row49[0,6] = 39200
$row49[0,7] = 38500
$row49[0,8] = 38900
$row49[0,9] = 39200
$ws.Range("D49:M49").Value2 = $row49

$row50 = New-Object 'object[,]' 1,10
$row50[0,0] = 0
$row50[0,1] = 0
$row50[0,2] = 0
$row50[0,3] = 0
$row50[0,4] = 0
$row50[0,5] = 0
$row50[0,6] = 0
$row50[0,7] = 0
$row50[0,8] = 0
$row50[0,9] = 0
$ws.Range("D50:M50").Value2 = $row50

$row51 = New-Object 'object[,]' 1,10
$row51[0,0] = 0
$row51[0,1] = 0
$row51[0,2] = 0
$row51[0,3] = 0
$row51[0,4] = 0
$row51[0,5] = 0
$row51[0,6] = 0
$row51[0,7] = 0
$row51[0,8] = 0
$row51[0,9] = 0
$ws.Range("D51:M51").Value2 = $row51

$row52 = New-Object 'object[,]' 1,10
$row52[0,0] = 20200
$row52[0,1] = 25200
$row52[0,2] = 17800
$row52[0,3] = 23100
$row52[0,4] = 78400
$row52[0,5] = 72400
$row52[0,6] = 67400
$row52[0,7] = 72000
$row52[0,8] = 73200
$row52[0,9] = 76600
$ws.Range("D52:M52").Value2 = $row52

$row53 = New-Object 'object[,]' 1,10
$row53[0,0] = 0
$row53[0,1] = 0
$row53[0,2] = 0
$row53[0,3] = 0
$row53[0,4] = 0
$row53[0,5] = 0
$row53[0,6] = 0
$row53[0,7] = 0
$row53[0,8] = 0
$row53[0,9] = 0
$ws.Range("D53:M53").Value2 = $row53

$row54 = New-Object 'object[,]' 1,10
$row54[0,0] = 5760600
$row54[0,1] = 6158700
$row54[0,2] = 6211800
$row54[0,3] = 6084700
$row54[0,4] = 5935200
$row54[0,5] = 5569100
$row54[0,6] = 2876500
$row54[0,7] = 2958000
$row54[0,8] = 2979800
$row54[0,9] = 3019500
$ws.Range("D54:M54").Value2 = $row54

$row57 = New-Object 'object[,]' 1,10
$row57[0,0] = 1009700
$row57[0,1] = 1108100
$row57[0,2] = 1050600
$row57[0,3] = 1005700
$row57[0,4] = 973400
$row57[0,5] = 800900
$row57[0,6] = 470400
$row57[0,7] = 503600
$row57[0,8] = 494600
$row57[0,9] = 383000
$ws.Range("D57:M57").Value2 = $row57

$row58 = New-Object 'object[,]' 1,10
$row58[0,0] = 32000
$row58[0,1] = 32000
$row58[0,2] = 180800
$row58[0,3] = 189400
$row58[0,4] = 590200
$row58[0,5] = 351000
$row58[0,6] = 91400
$row58[0,7] = 84400
$row58[0,8] = 84400
$row58[0,9] = 84400
$ws.Range("D58:M58").Value2 = $row58

$row59 = New-Object 'object[,]' 1,10
$row59[0,0] = 621800
$row59[0,1] = 865500
$row59[0,2] = 930000
$row59[0,3] = 943800
$row59[0,4] = 1108100
$row59[0,5] = 932300
$row59[0,6] = 314000
$row59[0,7] = 330900
$row59[0,8] = 356200
$row59[0,9] = 560300
$ws.Range("D59:M59").Value2 = $row59

$row60 = New-Object 'object[,]' 1,10
$row60[0,0] = 1663500
$row60[0,1] = 2005600
$row60[0,2] = 2161400
$row60[0,3] = 2138900
$row60[0,4] = 2671700
$row60[0,5] = 2084200
$row60[0,6] = 875800
$row60[0,7] = 918900
$row60[0,8] = 935200
$row60[0,9] = 1027700
$ws.Range("D60:M60").Value2 = $row60

$row61 = New-Object 'object[,]' 1,10
$row61[0,0] = 1751300
$row61[0,1] = 1830000
$row61[0,2] = 1861700
$row61[0,3] = 1770800
$row61[0,4] = 875400
$row61[0,5] = 1076800
$row61[0,6] = 731100
$row61[0,7] = 740500
$row61[0,8] = 748500
$row61[0,9] = 743300
$ws.Range("D61:M61").Value2 = $row61

$row62 = New-Object 'object[,]' 1,10
$row62[0,0] = 537700
$row62[0,1] = 426000
$row62[0,2] = 409400
$row62[0,3] = 380700
$row62[0,4] = 423900
$row62[0,5] = 624500
$row62[0,6] = 118200
$row62[0,7] = 116200
$row62[0,8] = 113600
$row62[0,9] = 101200
$ws.Range("D62:M62").Value2 = $row62

$row63 = New-Object 'object[,]' 1,10
$row63[0,0] = 0
$row63[0,1] = 0
$row63[0,2] = 0
$row63[0,3] = 0
$row63[0,4] = 0
$row63[0,5] = 0
$row63[0,6] = 0
$row63[0,7] = 0
$row63[0,8] = 0
$row63[0,9] = 0
$ws.Range("D63:M63").Value2 = $row63

$row64 = New-Object 'object[,]' 1,10
$row64[0,0] = 0
$row64[0,1] = 0
$row64[0,2] = 0
$row64[0,3] = 0
$row64[0,4] = 0
$row64[0,5] = 0
$row64[0,6] = 0
$row64[0,7] = 0
$row64[0,8] = 0
$row64[0,9] = 0
$ws.Range("D64:M64").Value2 = $row64

$row65 = New-Object 'object[,]' 1,10
$row65[0,0] = 0
$row65[0,1] = 0
$row65[0,2] = 0
$row65[0,3] = 0
$row65[0,4] = 0
$row65[0,5] = 0
$row65[0,6] = 0
$row65[0,7] = 0
$row65[0,8] = 0
$row65[0,9] = 0
$ws.Range("D65:M65").Value2 = $row65

$row66 = New-Object 'object[,]' 1,10
$row66[0,0] = 4128000
$row66[0,1] = 4437400
$row66[0,2] = 4609000
$row66[0,3] = 4466500
$row66[0,4] = 4284600
$row66[0,5] = 4088700
$row66[0,6] = 1905400
$row66[0,7] = 1960100
$row66[0,8] = 1987900
$row66[0,9] = 2071300
$ws.Range("D66:M66").Value2 = $row66

$row68 = New-Object 'object[,]' 1,10
$row68[0,0] = 0
$row68[0,1] = 0
$row68[0,2] = 0
$row68[0,3] = 0
$row68[0,4] = 0
$row68[0,5] = 0
$row68[0,6] = 0
$row68[0,7] = 0
$row68[0,8] = 0
$row68[0,9] = 0
$ws.Range("D68:M68").Value2 = $row68

$row69 = New-Object 'object[,]' 1,10
$row69[0,0] = 0
$row69[0,1] = 0
$row69[0,2] = 0
$row69[0,3] = 0
$row69[0,4] = 0
$row69[0,5] = 0
$row69[0,6] = 0
$row69[0,7] = 0
$row69[0,8] = 0
$row69[0,9] = 0
$ws.Range("D69:M69").Value2 = $row69

$row70 = New-Object 'object[,]' 1,10
$row70[0,0] = 0
$row70[0,1] = 0
$row70[0,2] = 0
$row70[0,3] = 0
$row70[0,4] = 0
$row70[0,5] = 0
$row70[0,6] = 0
$row70[0,7] = 0
$row70[0,8] = 0
$row70[0,9] = 0
$ws.Range("D70:M70").Value2 = $row70

$row71 = New-Object 'object[,]' 1,10
$row71[0,0] = 0
$row71[0,1] = 0
$row71[0,2] = 0
$row71[0,3] = 0
$row71[0,4] = 0
$row71[0,5] = 0
$row71[0,6] = 0
$row71[0,7] = 0
$row71[0,8] = 0
$row71[0,9] = 0
$ws.Range("D71:M71").Value2 = $row71

$row72 = New-Object 'object[,]' 1,10
$row72[0,0] = 981800
$row72[0,1] = 901500
$row72[0,2] = 742800
$row72[0,3] = 684300
$row72[0,4] = 767800
$row72[0,5] = 568600
$row72[0,6] = 476600
$row72[0,7] = 523900
$row72[0,8] = 522300
$row72[0,9] = 487500
$ws.Range("D72:M72").Value2 = $row72

$row73 = New-Object 'object[,]' 1,10
$row73[0,0] = 0
$row73[0,1] = 0
$row73[0,2] = 0
$row73[0,3] = 0
$row73[0,4] = 0
$row73[0,5] = 0
$row73[0,6] = 0
$row73[0,7] = 0
$row73[0,8] = 0
$row73[0,9] = 0
$ws.Range("D73:M73").Value2 = $row73

$row74 = New-Object 'object[,]' 1,10
$row74[0,0] = 0
$row74[0,1] = 0
$row74[0,2] = 0
$row74[0,3] = 0
$row74[0,4] = 0
$row74[0,5] = 0
$row74[0,6] = 0
$row74[0,7] = 0
$row74[0,8] = 0
$row74[0,9] = 0
$ws.Range("D74:M74").Value2 = $row74

$row75 = New-Object 'object[,]' 1,10
$row75[0,0] = 0
$row75[0,1] = 0
$row75[0,2] = 0
$row75[0,3] = 0
$row75[0,4] = 0
$row75[0,5] = 0
$row75[0,6] = 0
$row75[0,7] = 0
$row75[0,8] = 0
$row75[0,9] = 0
$ws.Range("D75:M75").Value2 = $row75

$row76 = New-Object 'object[,]' 1,10
$row76[0,0] = 1632600
$row76[0,1] = 1721300
$row76[0,2] = 1602800
$row76[0,3] = 1618200
$row76[0,4] = 1650600
$row76[0,5] = 1480400
$row76[0,6] = 971100
$row76[0,7] = 997900
$row76[0,8] = 991900
$row76[0,9] = 948200
$ws.Range("D76:M76").Value2 = $row76

$row77 = New-Object 'object[,]' 1,10
$row77[0,0] = 0
$row77[0,1] = 0
$row77[0,2] = 0
$row77[0,3] = 0
$row77[0,4] = 0
$row77[0,5] = 0
$row77[0,6] = 0
$row77[0,7] = 0
$row77[0,8] = 0
$row77[0,9] = 0
$ws.Range("D77:M77").Value2 = $row77

$row80 = New-Object 'object[,]' 1,10
$row80[0,0] = 43465
$row80[0,1] = 43373
$row80[0,2] = 43281
$row80[0,3] = 43190
$row80[0,4] = 43100
$row80[0,5] = 43008
$row80[0,6] = 42916
$row80[0,7] = 42825
$row80[0,8] = 42735
$row80[0,9] = 42643
$ws.Range("D80:M80").Value2 = $row80

$row81 = New-Object 'object[,]' 1,10
$row81[0,0] = 116100
$row81[0,1] = 179800
$row81[0,2] = 79100
$row81[0,3] = -34900
$row81[0,4] = 211100
$row81[0,5] = 104400
$row81[0,6] = -37900
$row81[0,7] = 11200
$row81[0,8] = 44200
$row81[0,9] = -161700
$ws.Range("D81:M81").Value2 = $row81

$row83 = New-Object 'object[,]' 1,10
$row83[0,0] = 53000
$row83[0,1] = 49200
$row83[0,2] = 49200
$row83[0,3] = 48000
$row83[0,4] = 47900
$row83[0,5] = 46900
$row83[0,6] = 29500
$row83[0,7] = 29000
$row83[0,8] = 29800
$row83[0,9] = 13200
$ws.Range("D83:M83").Value2 = $row83

$row84 = New-Object 'object[,]' 1,10
$row84[0,0] = 0
$row84[0,1] = 0
$row84[0,2] = 0
$row84[0,3] = 0
$row84[0,4] = 0
$row84[0,5] = 0
$row84[0,6] = 0
$row84[0,7] = 0
$row84[0,8] = 0
$row84[0,9] = 0
$ws.Range("D84:M84").Value2 = $row84

$row85 = New-Object 'object[,]' 1,10
$row85[0,0] = 0
$row85[0,1] = 0
$row85[0,2] = 0
$row85[0,3] = 0
$row85[0,4] = 0
$row85[0,5] = 0
$row85[0,6] = 0
$row85[0,7] = 0
$row85[0,8] = 0
$row85[0,9] = 0
$ws.Range("D85:M85").Value2 = $row85

$row86 = New-Object 'object[,]' 1,10
$row86[0,0] = 0
$row86[0,1] = 0
$row86[0,2] = 0
$row86[0,3] = 0
$row86[0,4] = 0
$row86[0,5] = 0
$row86[0,6] = 0
$row86[0,7] = 0
$row86[0,8] = 0
$row86[0,9] = 0
$ws.Range("D86:M86").Value2 = $row86

$row87 = New-Object 'object[,]' 1,10
$row87[0,0] = 0
$row87[0,1] = 0
$row87[0,2] = 0
$row87[0,3] = 0
$row87[0,4] = 0
$row87[0,5] = 0
$row87[0,6] = 0
$row87[0,7] = 0
$row87[0,8] = 0
$row87[0,9] = 0
$ws.Range("D87:M87").Value2 = $row87

$row88 = New-Object 'object[,]' 1,10
$row88[0,0] = 0
$row88[0,1] = 0
$row88[0,2] = 0
$row88[0,3] = 0
$row88[0,4] = 0
$row88[0,5] = 0
$row88[0,6] = 0
$row88[0,7] = 0
$row88[0,8] = 0
$row88[0,9] = 0
$ws.Range("D88:M88").Value2 = $row88

$row89 = New-Object 'object[,]' 1,10
$row89[0,0] = 359100
$row89[0,1] = 337600
$row89[0,2] = 54300
$row89[0,3] = -190700
$row89[0,4] = 248800
$row89[0,5] = 119600
$row89[0,6] = 5800
$row89[0,7] = -42100
$row89[0,8] = 146700
$row89[0,9] = -35100
$ws.Range("D89:M89").Value2 = $row89

$row91 = New-Object 'object[,]' 1,10
$row91[0,0] = -94000
$row91[0,1] = -85500
$row91[0,2] = -71500
$row91[0,3] = -71000
$row91[0,4] = -63600
$row91[0,5] = -73400
$row91[0,6] = -16000
$row91[0,7] = -38000
$row91[0,8] = -18100
$row91[0,9] = 19500
$ws.Range("D91:M91").Value2 = $row91

$row92 = New-Object 'object[,]' 1,10
$row92[0,0] = 0
$row92[0,1] = 0
$row92[0,2] = 0
$row92[0,3] = 0
$row92[0,4] = 0
$row92[0,5] = 0
$row92[0,6] = 0
$row92[0,7] = 0
$row92[0,8] = 0
$row92[0,9] = 0
$ws.Range("D92:M92").Value2 = $row92

$row93 = New-Object 'object[,]' 1,10
$row93[0,0] = 0
$row93[0,1] = 0
$row93[0,2] = 0
$row93[0,3] = 0
$row93[0,4] = 0
$row93[0,5] = 0
$row93[0,6] = 0
$row93[0,7] = 0
$row93[0,8] = 0
$row93[0,9] = 0
$ws.Range("D93:M93").Value2 = $row93

$row94 = New-Object 'object[,]' 1,10
$row94[0,0] = -88100
$row94[0,1] = -51000
$row94[0,2] = 40400
$row94[0,3] = -26600
$row94[0,4] = -68600
$row94[0,5] = 138300
$row94[0,6] = -16500
$row94[0,7] = -15600
$row94[0,8] = 278100
$row94[0,9] = -38200
$ws.Range("D94:M94").Value2 = $row94

$row96 = New-Object 'object[,]' 1,10
$row96[0,0] = -21300
$row96[0,1] = -21000
$row96[0,2] = -20800
$row96[0,3] = -17000
$row96[0,4] = -12700
$row96[0,5] = -12300
$row96[0,6] = -9400
$row96[0,7] = -9600
$row96[0,8] = -9400
$row96[0,9] = -9200
$ws.Range("D96:M96").Value2 = $row96

$row97 = New-Object 'object[,]' 1,10
$row97[0,0] = 0
$row97[0,1] = 0
$row97[0,2] = 0
$row97[0,3] = 0
$row97[0,4] = 0
$row97[0,5] = 0
$row97[0,6] = 0
$row97[0,7] = 0
$row97[0,8] = 0
$row97[0,9] = 0
$ws.Range("D97:M97").Value2 = $row97

$row98 = New-Object 'object[,]' 1,10
$row98[0,0] = 0
$row98[0,1] = 0
$row98[0,2] = 0
$row98[0,3] = 0
$row98[0,4] = 0
$row98[0,5] = 0
$row98[0,6] = 0
$row98[0,7] = 0
$row98[0,8] = 0
$row98[0,9] = 0
$ws.Range("D98:M98").Value2 = $row98

$row99 = New-Object 'object[,]' 1,10
$row99[0,0] = 0
$row99[0,1] = 0
$row99[0,2] = 0
$row99[0,3] = 0
$row99[0,4] = 0
$row99[0,5] = 0
$row99[0,6] = 0
$row99[0,7] = 0
$row99[0,8] = 0
$row99[0,9] = 0
$ws.Range("D99:M99").Value2 = $row99

$row100 = New-Object 'object[,]' 1,10
$row100[0,0] = -300800
$row100[0,1] = -310300
$row100[0,2] = 20100
$row100[0,3] = 293400
$row100[0,4] = -74800
$row100[0,5] = 8100
$row100[0,6] = -2900
$row100[0,7] = -35000
$row100[0,8] = -65900
$row100[0,9] = 26500
$ws.Range("D100:M100").Value2 = $row100

$row101 = New-Object 'object[,]' 1,10
$row101[0,0] = 0
$row101[0,1] = 0
$row101[0,2] = 0
$row101[0,3] = 0
$row101[0,4] = 0
$row101[0,5] = 0
$row101[0,6] = 0
$row101[0,7] = 0
$row101[0,8] = 0
$row101[0,9] = 0
$ws.Range("D101:M101").Value2 = $row101

$row102 = New-Object 'object[,]' 1,10
$row102[0,0] = -29800
$row102[0,1] = -23700
$row102[0,2] = 114800
$row102[0,3] = 76100
$row102[0,4] = 103900
$row102[0,5] = 265700
$row102[0,6] = -19100
$row102[0,7] = -97800
$row102[0,8] = 358900
$row102[0,9] = -46800
$ws.Range("D102:M102").Value2 = $row102

